$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://www.youtube.com/watch?v=FSVtc3pDO8w&t=&ab_channel=DeitiesofDeath"

$target = $ws.Range("A19")

$ws.Hyperlinks.Add($target, $url, "", $url, $url) | Out-Null

$target.Select()
